$d = $word.ActiveDocument

# "As a user I want to bookmark my favorite jobs + save in local storage"
# becomes
# "As a user I want to save my favorite jobs and when I reload page still see my saved job"

# Step 1: change "bookmark" -> "save" in the lead-in sentence.
$d.Content.Find.Execute("- As a user I want to bookmark my ", $true, $false, $false, $false, $false, $true, 1, $false, "- As a user I want to save my ", 2)

# Step 2: replace the tail (after "favorite") describing the old "+ save in
# local storage" note (including the stray grammar-check markers around
# "jobs") with the new sentence about reloading the page.
$d.Content.Find.Execute(" jobs  + save in local storage", $true, $false, $false, $false, $false, $true, 1, $false, " jobs and when I reload page still see my saved job", 2)
